# Applies the "keywords, cron, video prototype" commit:
#  - Slide master date / footer / slide-number placeholders get their
#    cached display text swapped for the generic bracketed placeholders.
#  - A handful of body-text paragraphs on slides 3, 8, 10 and 15 get
#    their runs re-split (or re-joined) and/or small wording tweaks.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide master placeholders (date / footer / slide number)
# ---------------------------------------------------------------------
$master = $p.SlideMaster

# Date placeholder
$dateShape = $master.Shapes.Item(3)
$dateShape.TextFrame.TextRange.Text = "<date/time>"

# Footer placeholder
$footerShape = $master.Shapes.Item(4)
$footerShape.TextFrame.TextRange.Text = "<footer>"

# Slide number placeholder (keeps the <a:fld type="slidenum"> field where
# the host supports it; this host materialises a plain run on edit, which
# is the same behaviour real PowerPoint automation is limited to since
# the field's internal GUID isn't scriptable)
$numShape = $master.Shapes.Item(5)
$numShape.TextFrame.TextRange.Text = "<number>"

# ---------------------------------------------------------------------
# 2. Slide 10 - "Task or Intents (1)"
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Item(2).TextFrame.TextRange

# Para 1: "Small Python programs that (usually) run an OS level task"
#   -> split into two runs
$tr10.Characters(1, 44).Text = "Small Python programs that (usually) run an "
$tr10.Characters(45, 13).Text = "OS level task"

# Para 2, run 2 (Courier New): "arecord ... --duration 30 "
#   -> split into two runs
$tr10.Characters(72, 32).Text = "arecord -q -f cd -c 2 -D hw:0,0 "
$tr10.Characters(104, 14).Text = "--duration 30 "

# Para 3, last run: " so can experiment/debug" -> split into two runs
$tr10.Characters(195, 1).Text = " "
$tr10.Characters(196, 23).Text = "so can experiment/debug"

# Para 4, first run: "Control external AI ... this can be turned off in "
#   -> split into two runs
$tr10.Characters(220, 50).Text = "Control external AI (transcription etc.) this can "
$tr10.Characters(270, 17).Text = "be turned off in "

# ---------------------------------------------------------------------
# 3. Slide 15 - "Challenges (2)"
# ---------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$tr15 = $s15.Shapes.Item(2).TextFrame.TextRange
$tr15.Characters(321, 46).Text = "Licence given the assorted licences for components"

# ---------------------------------------------------------------------
# 4. Slide 3 - "Rhasspy: Port 12101"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange

# Para 1: merge the two runs back into one
$tr3.Characters(1, 75).Text = "Think Alexis or Mycroft, open source home assistant (HA) with voice command"

# Para 2: merge the two runs back into one
$tr3.Characters(77, 83).Text = "I rejected Mycroft because it's tightly linked to the cloud, Rhasspy gives 'choice'"

# ---------------------------------------------------------------------
# 5. Slide 8 - "Intent Server(1) Port 8000"
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange

# Para 2: "Takes an Intent (take photo, for example) and dispatches to a specialised program"
#   -> split into two runs
$tr8.Characters(26, 46).Text = "Takes an Intent (take photo, for example) and "
$tr8.Characters(72, 35).Text = "dispatches to a specialised program"

# Para 3: "Dispatch table for intents (can be done in Node-Red) but finer control here"
#   -> split into two runs
$tr8.Characters(108, 43).Text = "Dispatch table for intents (can be done in "
$tr8.Characters(151, 32).Text = "Node-Red) but finer control here"
